{"js": "// Apply strikethrough formatting to the six heading paragraphs that were\n// struck through in the target revision: \"Summary\", \"Appendices\", \"Code\",\n// \"Dataset\", \"Model Accuracy Visualization\" and \"Model Feature Importance\n// Visualization\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Summary\",\n  \"Appendices\",\n  \"Code\",\n  \"Dataset\",\n  \"Model Accuracy Visualization\",\n  \"Model Feature Importance Visualization\",\n]);\n\nfor (const paragraph of paragraphs.items) {\n  const text = (paragraph.text || \"\").trim();\n  if (targets.has(text)) {\n    paragraph.font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the six heading paragraphs that were\n# struck through in the target revision: \"Summary\", \"Appendices\", \"Code\",\n# \"Dataset\", \"Model Accuracy Visualization\" and \"Model Feature Importance\n# Visualization\".\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Summary\",\n    \"Appendices\",\n    \"Code\",\n    \"Dataset\",\n    \"Model Accuracy Visualization\",\n    \"Model Feature Importance Visualization\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\")\n    if ($targets -contains $text) {\n        $p.Range.Font.StrikeThrough = 1\n    }\n}\n"}
